# Handles float input without breaking stuff
#
# Rebuilds the marksheet's score summary (rows 10-12), removes the
# second and third "Student Ans / Correct Ans" blocks (columns D:H
# below row 15), and repurposes D16:E18 plus a handful of rows in the
# first block (A:B) with the data that used to live in those removed
# blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Score summary block (rows 10-12)
# ---------------------------------------------------------------------

# Row 10: No. / Right / Wrong / Not Attempt / Max
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

# Row 11: Marking (C11 used to be a literal text "-1", now a real
# number so it keeps behaving under SUM()/etc.)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Total
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "37/112"

# A10:A12 keep their text but pick up the "mtitleStyle" formatting
# already used by A9 (same column/row style, s="4") instead of being
# unstyled.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Drop the third "Student Ans / Correct Ans" block (columns G:H)
#    entirely. Clearing the whole columns removes the cells from the
#    sheet outright (instead of just blanking them) and shrinks the
#    used range from H40 down to E40.
# ---------------------------------------------------------------------
$ws.Range("G1:H1048576").Clear()

# ---------------------------------------------------------------------
# 3) Drop the now-unused tail of the second block (columns D:E, rows
#    19-40) - only rows 16-18 of that block survive, repurposed below.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 4) Repurpose D16:E18 (former second "Student Ans/Correct Ans" block)
#    with the answers that used to sit in the deleted G:H block, using
#    existing correctStyle/incorrectStyle/absoluteStyle cells as
#    format donors so no new style records get created.
# ---------------------------------------------------------------------

# D16 -> correctStyle, "Option A"
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D16").Value = "Option A"

# D17 -> incorrectStyle, "Option A"
$ws.Range("C10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D17").Value = "Option A"

# D18 -> correctStyle, "Option D"
$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D18").Value = "Option D"

# E16:E18 -> absoluteStyle (same style already used further up in
# E12), "Option A" / "Option C" / "Option D"
$ws.Range("E12").Copy()
$ws.Range("E16:E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E16").Value = "Option A"
$ws.Range("E17").Value = "Option C"
$ws.Range("E18").Value = "Option D"

# ---------------------------------------------------------------------
# 5) A handful of rows in the first "Student Ans / Correct Ans" block
#    switch column A from an unstyled blank placeholder to a
#    correctStyle/incorrectStyle cell carrying the actual answer, and
#    gain a matching "Correct Ans" entry (absoluteStyle) in column B.
# ---------------------------------------------------------------------

function Set-AnswerPair([int]$row, [string]$styleSourceCell, [string]$aText, [string]$bText) {
    $localWs = $wb.ActiveSheet

    $localWs.Range($styleSourceCell).Copy()
    $localWs.Range("A" + $row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $localWs.Range("A" + $row).Value = $aText

    $localWs.Range("E16").Copy()
    $localWs.Range("B" + $row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $localWs.Range("B" + $row).Value = $bText
}

# correctStyle (s=5) rows - use B10 (correctStyle) as format donor
Set-AnswerPair 21 "B10" "Option C" "Option C"
Set-AnswerPair 22 "B10" "Option D" "Option D"
Set-AnswerPair 23 "B10" "Option D" "Option D"
Set-AnswerPair 29 "B10" "Option D" "Option D"
Set-AnswerPair 32 "B10" "Option C" "Option C"
Set-AnswerPair 37 "B10" "Option A" "Option A"
Set-AnswerPair 38 "B10" "Option A" "Option A"
Set-AnswerPair 39 "B10" "Option D" "Option D"

# incorrectStyle (s=6) rows - use C10 (incorrectStyle) as format donor
Set-AnswerPair 24 "C10" "Option D" "Option A"
Set-AnswerPair 28 "C10" "Option B" "Option D"
